$wb = $excel.ActiveWorkbook

# Sheet1: Balance for Sample1 (F3) increased from 2300 to 2500
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("F3").Value = 2500

# Sheet2: Overdraft/Balance figures updated for ABC company
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("F3").Value = 2000
$ws2.Range("G3").Value = 0

# Update the active selection on Sheet2 to G3, and make Sheet2 the active sheet
$ws2.Activate()
$ws2.Range("G3").Select()
